$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank -> "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Old row 10 was Contact / "No display for ContactDetail" -> becomes Jurisdiction / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Old row 11 was a duplicate Contact row - remove it entirely, shifting everything below up by one
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

# Row 2 (the root Extension element) Short/Definition text update
$elem.Range("K2").Value = "Employee Status"
$elem.Range("L2").Value = "Status of the employee based on one or more code systems. Example codes include HIPAA (HipaaEmployeeStatusCodeSystem), Payer (WhPayerEmployeeStatusCodeSystem) or customer-specific codes."
